$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Materias primas" (raw materials) column values - reordering ingredient lists
$ws.Range("C2").Value = "harina,huevos,vainilla,leche,"
$ws.Range("C4").Value = "huevos,vainilla,harina,"
$ws.Range("C5").Value = "harina,huevos,"
$ws.Range("C6").Value = "harina, huevos, limon,merengue,crema,"
